$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.347.17"
$ws.Range("E2").Value = "  +0.64%  "

# Row 3
$ws.Range("D3").Value = "1.621.18"
$ws.Range("E3").Value = "  +1.19%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'212.24"
$ws.Range("E5").Value = "  +0.18%  "

# Row 6
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("E7").Value = "  +0.20%  "

# Row 8
$ws.Range("E8").Value = "  +0.20%  "

# Row 9
$ws.Range("D9").Value = "'0.0616"
$ws.Range("E9").Value = "  +0.34%  "

# Row 10
$ws.Range("D10").Value = "'18.75"
$ws.Range("E10").Value = "  +3.47%  "

# Row 11
$ws.Range("E11").Value = "  +0.31%  "

# Row 12
$ws.Range("D12").Value = "1.850.40"
$ws.Range("E12").Value = "  +1.46%  "

# Row 13
$ws.Range("D13").Value = "1.613.03"
$ws.Range("E13").Value = "  +0.61%  "

# Row 14
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "  +0.27%  "

# Row 15
$ws.Range("E15").Value = "  +0.63%  "

# Row 16
$ws.Range("D16").Value = "26.361.64"
$ws.Range("E16").Value = "  +0.76%  "

# Row 17
$ws.Range("D17").Value = "'62.50"
$ws.Range("E17").Value = "  +2.56%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  +0.20%  "

# Row 19
$ws.Range("E19").Value = "  -0.16%  "

# Row 20
$ws.Range("D20").Value = "'202.26"
$ws.Range("E20").Value = "  -1.04%  "

# Row 21
$ws.Range("D21").Value = "'4.26"
$ws.Range("E21").Value = "  -0.06%  "

# Row 22
$ws.Range("D22").Value = "'9.28"
$ws.Range("E22").Value = "  +0.14%  "

# Row 23
$ws.Range("D23").Value = "'6.05"
$ws.Range("E23").Value = "  +0.40%  "

# Row 24
$ws.Range("E24").Value = "  -3.86%  "

# Row 25
$ws.Range("D25").Value = "'144.44"
$ws.Range("E25").Value = "  +0.06%  "

# Row 26
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("E27").Value = "  -1.96%  "

# Row 28
$ws.Range("E28").Value = "  -0.02%  "

# Row 29
$ws.Range("E29").Value = "  +1.15%  "

# Row 30
$ws.Range("E30").Value = "  +6.30%  "

# Row 31
$ws.Range("E31").Value = "  +0.15%  "

# Row 32
$ws.Range("E32").Value = "  +1.39%  "

# Row 33
$ws.Range("E33").Value = "  +0.41%  "

# Row 34
$ws.Range("E34").Value = "  +0.39%  "

# Row 35
$ws.Range("E35").Value = "  +2.12%  "

# Row 36
$ws.Range("D36").Value = "1.162.80"
$ws.Range("E36").Value = "  +1.95%  "

# Row 37
$ws.Range("E37").Value = "  +0.05%  "

# Row 38
$ws.Range("D38").Value = "'0.801"
$ws.Range("E38").Value = "  +1.83%  "

# Row 39
$ws.Range("E39").Value = "  -0.10%  "

# Row 40
$ws.Range("E40").Value = "  -0.05%  "

# Row 41
$ws.Range("D41").Value = "'0.495"
$ws.Range("E41").Value = "  +0.15%  "

# Row 42
$ws.Range("D42").Value = "'5.40"
$ws.Range("E42").Value = "  +4.17%  "

# Row 43
$ws.Range("E43").Value = "  +0.00%  "

# Row 44
$ws.Range("D44").Value = "1.761.81"
$ws.Range("E44").Value = "  +1.40%  "

# Row 45
$ws.Range("D45").Value = "'92.19"
$ws.Range("E45").Value = "  +0.08%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0104"
$ws.Range("E46").Value = "  +9.73%  "

# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.52"
$ws.Range("E47").Value = "  +0.90%  "

# Row 48
$ws.Range("D48").Value = "'53.82"
$ws.Range("E48").Value = "  -0.46%  "

# Row 49
$ws.Range("D49").Value = "'0.0507"
$ws.Range("E49").Value = "  +0.09%  "

# Row 50
$ws.Range("E50").Value = "  +0.93%  "

# Row 51
$ws.Range("E51").Value = "  -0.27%  "
